# The commit swaps the contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
#   - theme1.xml goes from "Office Theme" (default Office colours) to "Integral"
#     (the colour scheme that theme2.xml/the slide master already used).
#   - theme2.xml (the theme actually bound to the presentation's one-and-only
#     slide master / design) goes from "Integral" back to the stock
#     "Office Theme" colour values.
#
# This PowerPoint COM host only exposes one editable colour theme - the one
# attached to the (single) slide master/design - via
# Slide/SlideMaster.ThemeColorScheme, and only the 12 theme colour slots are
# settable (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). So we drive the
# presentation's design back to the "Office Theme" palette through that API,
# matching the colours that ship in the target ppt/theme/theme2.xml.

$p = $ppt.ActivePresentation

function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgr = $r + ($g * 256) + ($b * 65536)
    $scheme.Colors($index).RGB = $bgr
}

$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Office Theme colour scheme (the values that used to live in theme1.xml,
# now destined for theme2.xml):
Set-ThemeColor $colorScheme 1  "000000"   # dk1
Set-ThemeColor $colorScheme 2  "FFFFFF"   # lt1
Set-ThemeColor $colorScheme 3  "44546A"   # dk2
Set-ThemeColor $colorScheme 4  "E7E6E6"   # lt2
Set-ThemeColor $colorScheme 5  "5B9BD5"   # accent1
Set-ThemeColor $colorScheme 6  "ED7D31"   # accent2
Set-ThemeColor $colorScheme 7  "A5A5A5"   # accent3
Set-ThemeColor $colorScheme 8  "FFC000"   # accent4
Set-ThemeColor $colorScheme 9  "4472C4"   # accent5
Set-ThemeColor $colorScheme 10 "70AD47"   # accent6
Set-ThemeColor $colorScheme 11 "0563C1"   # hlink
Set-ThemeColor $colorScheme 12 "954F72"   # folHlink
